$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 9.800000000000001
$ws.Range("G2").Value = 16.5
$ws.Range("H2").Value = 1.29
$ws.Range("I2").Value = 1.35
$ws.Range("J2").Value = 5.8
$ws.Range("K2").Value = 6.8
$ws.Range("P2").Value = 2.14
$ws.Range("Q2").Value = 1.51
$ws.Range("BH2").Value = "2026-02-24 10:13:34"

# Row 3 (only timestamp)
$ws.Range("BH3").Value = "2026-02-24 10:13:34"

# Row 4
$ws.Range("G4").Value = 3.85
$ws.Range("BH4").Value = "2026-02-24 10:13:34"

# Row 5
$ws.Range("F5").Value = 1.85
$ws.Range("G5").Value = 1.93
$ws.Range("H5").Value = 5.4
$ws.Range("K5").Value = 3.6
$ws.Range("P5").Value = 1.63
$ws.Range("Q5").Value = 2.4
$ws.Range("BH5").Value = "2026-02-24 10:13:34"

# Row 6 (only timestamp)
$ws.Range("BH6").Value = "2026-02-24 10:13:34"

# Row 7
$ws.Range("F7").Value = 2.12
$ws.Range("G7").Value = 2.26
$ws.Range("H7").Value = 3.65
$ws.Range("K7").Value = 3.7
$ws.Range("Q7").Value = 2.08
$ws.Range("BH7").Value = "2026-02-24 10:13:34"

# Row 8
$ws.Range("G8").Value = 4.4
$ws.Range("H8").Value = 2.36
$ws.Range("I8").Value = 2.56
$ws.Range("J8").Value = 2.7
$ws.Range("K8").Value = 2.96
$ws.Range("BH8").Value = "2026-02-24 10:13:34"

# Row 9
$ws.Range("F9").Value = 1.5
$ws.Range("G9").Value = 1.56
$ws.Range("H9").Value = 8.800000000000001
$ws.Range("I9").Value = 10
$ws.Range("K9").Value = 4.4
$ws.Range("P9").Value = 1.64
$ws.Range("Q9").Value = 2.28
$ws.Range("BH9").Value = "2026-02-24 10:13:34"

# Row 10 (only timestamp)
$ws.Range("BH10").Value = "2026-02-24 10:13:34"
